# MasterAchievement.xlsx: introduce a "None" enum entry at the top of the
# achievement list (id 90000) and shift the existing GoodEnding / BadEnding /
# Lockmaster rows down by one, which in effect inserts a brand-new row above
# the previously-last two rows (ClearGameOnce / AllAchievements) while the
# header note block (K1:S2) stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for one more data row by inserting a blank row right above the
# "ClearGameOnce" entry (old row 5). This pushes the old row 5
# (ClearGameOnce, custom "Good" style) and row 6 (AllAchievements) down to
# rows 6 and 7, carrying their formatting with them, while rows 1-4 (and the
# K:S note cells on rows 1-2) are left completely alone by this operation.
$ws.Rows("5:5").Insert()

# Row 2 becomes the new "None" entry (id 90000 keeps its value, but now
# represents "no achievement"), so it loses its steamID key text.
$ws.Range("B2").Value = "None"
$ws.Range("E2").ClearContents()

# Row 3 takes over what used to be the GoodEnding row's text values.
$ws.Range("B3").Value = "GoodEnding"
$ws.Range("E3").Value = "goodEnding_90000"

# Row 4 takes over what used to be the BadEnding row's text values, and its
# id is renumbered to 90002.
$ws.Range("A4").Value = 90002
$ws.Range("B4").Value = "BadEnding"
$ws.Range("E4").Value = "badEnding_90001"
$ws.Range("G4").Value = -1
$ws.Range("H4").Value = 0

# Row 5 (the freshly inserted, previously blank row) gets the data that used
# to live in row 4 (Lockmaster).
$ws.Range("A5").Value = 91000
$ws.Range("B5").Value = "Lockmaster"
$ws.Range("C5").Value = -1
$ws.Range("D5").Value = -1
$ws.Range("E5").Value = "lockMaster_91000"
$ws.Range("F5").Value = $false
$ws.Range("G5").Value = 10201
$ws.Range("H5").Value = 1

# Match the author's final cursor position.
$ws.Range("J7").Select()
